$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header C1: "timestamp" -> "time"
$ws.Range("C1").Value = "time"

# Update row 2 data
$ws.Range("A2").Value = "小A"
$ws.Range("B2").Value = "今天排班很順利"
$ws.Range("C2").Value = "2025-07-17 14:00"

# Add new row 3 data
$ws.Range("A3").Value = "小B"
$ws.Range("B3").Value = "建議明天增加會議時間"
$ws.Range("C3").Value = "2025-07-17 15:30"
